# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap "Santa Lucia" (row 202) and "Timor Oriental" (row 203) country names
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# Update "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 16:56"

# Update numeric data cells per country row

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 5100910
$ws.Range("C4").Value = 5386
$ws.Range("D4").Value = 2618191
$ws.Range("E4").Value = 2318523
$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 164196

# Row 6 - Rusia
$ws.Range("B6").Value = 2114140
$ws.Range("C6").Value = 27276
$ws.Range("E6").Value = 628136
$ws.Range("G6").Value = 243
$ws.Range("H6").Value = 42821

# Row 15
$ws.Range("B15").Value = 309763
$ws.Range("C15").Value = 758

# Row 21
$ws.Range("D21").Value = 108242
$ws.Range("E21").Value = 122985
$ws.Range("G21").Value = 39
$ws.Range("H21").Value = 4450

# Row 22
$ws.Range("B22").Value = 216562
$ws.Range("C22").Value = 247
$ws.Range("E22").Value = 9905
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 9257

# Row 46
$ws.Range("B46").Value = 56189
$ws.Range("C46").Value = 919
$ws.Range("D46").Value = 44072
$ws.Range("E46").Value = 9920
$ws.Range("G46").Value = 29
$ws.Range("H46").Value = 2197

# Row 47
$ws.Range("D47").Value = 48583
$ws.Range("E47").Value = 6319

# Row 48
$ws.Range("B48").Value = 52537
$ws.Range("C48").Value = 186
$ws.Range("D48").Value = 38364
$ws.Range("E48").Value = 12423
$ws.Range("G48").Value = 4
$ws.Range("H48").Value = 1750

# Row 64
$ws.Range("B64").Value = 27443
$ws.Range("C64").Value = 453
$ws.Range("D64").Value = 19100
$ws.Range("E64").Value = 7502
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 841

# Row 73
$ws.Range("B73").Value = 19978
$ws.Range("C73").Value = 434
$ws.Range("D73").Value = 9515
$ws.Range("E73").Value = 9927

# Row 86
$ws.Range("B86").Value = 9568
$ws.Range("C86").Value = 17
$ws.Range("E86").Value = 455

# Row 126
$ws.Range("B126").Value = 2247
$ws.Range("C126").Value = 14
$ws.Range("E126").Value = 194
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 126

# Row 170
$ws.Range("D170").Value = 311
$ws.Range("E170").Value = 42
